$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date for the
# 8802293c-4c98-4afd-912e-c4e51c8b5e3d.md row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-06 19:01:03"

# Sheet "zh-cn": Correspond Handoff/Handback DateTime for the
# 8802293c-4c98-4afd-912e-c4e51c8b5e3d row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-06 19:00:55"
$wsZhCn.Range("K3").Value = "2016-09-06 19:01:37"

# Sheet "de-de": Correspond Handoff/Handback DateTime for the
# 8802293c-4c98-4afd-912e-c4e51c8b5e3d row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-06 19:01:03"
$wsDeDe.Range("K3").Value = "2016-09-06 19:01:46"
